# Commit: "nhcuong: update path separator"
#
# The sheet "Phieu Danh Gia" (Worksheets(1)) has a column F ("GHI CHU" /
# notes) that held explanatory strings about file/folder path handling
# (e.g. "khong cho phep thuc hien tren thu muc", "delete multi files",
# "move multi folders", ...). Those notes are removed from rows 10-15 and
# 17-20 - the cells become blank. Once those shared strings are no longer
# referenced anywhere in the workbook they drop out of the shared string
# table on save, which is what re-numbers the remaining shared-string
# indices used elsewhere on the sheet (teacher name cells, etc.) - that
# renumbering happens automatically when we clear the cells below, we
# don't need to touch those cells ourselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Clear the stale "GHI CHU" notes in column F (keeps the cell's style,
# only drops the text value - matches the target cells which stay
# `s="14"`/`s="20"` with no `t="s"`/`<v>` child).
$ws.Range("F10:F15").ClearContents()
$ws.Range("F17:F20").ClearContents()

# Update the view state: scroll back up near the top of the sheet and
# leave the note column selected (F10:F21) as the active selection.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F10:F21").Select()
